$d = $word.ActiveDocument

# 1. Update the generated timestamp (paragraph 2, single run - safe via Find)
$oldDate = 'Generated: 2025-09-01 11:36:30'
$newDate = 'Generated: 2025-09-04 10:57:32'
$null = $d.Content.Find.Execute($oldDate, $true, $false, $false, $false, $false, $true, 1, $false, $newDate, 2)

# 2. Replace the entire 'Product Overview' body paragraph content (paragraph 5)
#    Using InsertXML on the full paragraph text range keeps xml:space="preserve"
#    attributes exact (as crafted below) and leaves sibling paragraphs untouched.
$p5Fragment = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>## Product Overview</w:t><w:br/><w:br/><w:t>This product is an in-house credit card core system designed specifically for the U.S. financial services industry within the credit card sector. It serves financial institutions aiming to transition away from legacy third-party platforms to gain full ownership and control over the entire credit card lifecycle. Delivered as an integrated, real-time platform, it supports issuance of both virtual and physical credit cards, as well as comprehensive account management including credit lines, authorizations, settlements, billing, payments, interest calculations, fees, rewards, disputes, delinquency workflows, and charged-off account management. The system’s capability to handle end-to-end credit card processes makes it a central operational backbone for credit issuance and servicing.</w:t><w:br/><w:br/><w:t xml:space="preserve">Strategically, this core system addresses the need for modernization and operational independence by embedding critical financial functions that were previously managed externally. By internalizing these capabilities, the product enables faster innovation, improved data integrity, and enhanced customer experience. It supports regulatory compliance relevant to the U.S. credit card market and positions the enterprise to adapt seamlessly to evolving business needs and competitive pressures in a dynamic financial services environment. This shift to a modern core ultimately drives greater operational efficiency and risk control across the credit portfolio. </w:t><w:br/><w:br/><w:t xml:space="preserve">### References  </w:t><w:br/><w:t>No external sources were used in the generation of this overview.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p5 = $d.Paragraphs(5)
$p5Range = $d.Range($p5.Range.Start, $p5.Range.End - 1)
$p5Range.InsertXML($p5Fragment)

# 3. Replace the entire 'Feature Overview' body paragraph content (paragraph 7)
$p7Fragment = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Feature Name: Create a Frontbook Charged-Off Account</w:t><w:br/><w:br/><w:t xml:space="preserve">Feature Overview:  </w:t><w:br/><w:t>The Create a Frontbook Charged-Off Account feature is a specialized subset of the broader account creation capability within the in-house credit card core system. It enables comprehensive lifecycle management of accounts that have transitioned into charged-off status due to triggers such as prolonged delinquency, customer death, or bankruptcy. This feature supports real-time processing and integration with credit line management, authorizations, settlements, billing, payments, interest calculations, fees, rewards, disputes, and delinquency workflows, ensuring seamless handling of charged-off accounts within the trade credit ecosystem.</w:t><w:br/><w:br/><w:t xml:space="preserve">This feature includes the creation, status transition, and ongoing management of charged-off accounts but excludes the initial issuance of accounts in good standing or unrelated account servicing functions. It integrates with core systems responsible for transaction posting, interest calculation, account updates, and dispute management, leveraging data inputs from daily transaction files and account status triggers. Critical constraints include compliance with financial regulations governing charged-off accounts, real-time data accuracy, and secure handling of sensitive customer information. Strategically, this feature supports the organization’s goal of full ownership over credit issuance and servicing by replacing legacy third-party systems with a modern, flexible core that enhances operational control, risk management, and customer lifecycle visibility. </w:t><w:br/><w:br/><w:t>### References</w:t><w:br/><w:t>No external sources were used.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p7 = $d.Paragraphs(7)
$p7Range = $d.Range($p7.Range.Start, $p7.Range.End - 1)
$p7Range.InsertXML($p7Fragment)

Write-Host "Edit complete"
